# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to reflect newly generated output.

$wb = $excel.ActiveWorkbook

# Updates shared by both sheets (row -> new value for column F)
$commonUpdates = @{
    5  = 129
    6  = 201
    7  = 1687
    8  = 1632
    9  = 466
    15 = 232
    18 = 8
    21 = 52
    24 = 210
    25 = 105
    26 = 32
    27 = 9
    29 = 258
    30 = 2168
    34 = 329
    39 = 413
    40 = 517
}

# Sheet "展览" -> F3 becomes 3192
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(3, 6).Value = 3192
foreach ($row in $commonUpdates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $commonUpdates[$row]
}

# Sheet "全部类型" -> F3 becomes 3193
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(3, 6).Value = 3193
foreach ($row in $commonUpdates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $commonUpdates[$row]
}
